# Insert new leading columns into the "association" sheets that gained
# part_of / happens_during slots.
#
#   BiologicalProcessAssociation:          + happens_during, part_of  (A1, B1)
#   CellularAnatomicalEntityAssociation:   + part_of                 (A1)
#   CellTypeAssociation:                   + part_of                 (A1)
#   GrossAnatomyAssociation:               + part_of                 (A1)
#
# In every case the new column(s) are inserted before the existing
# "term" column, pushing term/type/evidence/provenances to the right.

$wb = $excel.ActiveWorkbook
$xlShiftToRight = -4161

# --- BiologicalProcessAssociation: insert two columns (happens_during, part_of) ---
$ws = $wb.Worksheets.Item("BiologicalProcessAssociation")
$ws.Range("A1:B1").Insert($xlShiftToRight)
$ws.Range("A1").Value = "happens_during"
$ws.Range("B1").Value = "part_of"

# --- CellularAnatomicalEntityAssociation: insert one column (part_of) ---
$ws = $wb.Worksheets.Item("CellularAnatomicalEntityAssociation")
$ws.Range("A1").Insert($xlShiftToRight)
$ws.Range("A1").Value = "part_of"

# --- CellTypeAssociation: insert one column (part_of) ---
$ws = $wb.Worksheets.Item("CellTypeAssociation")
$ws.Range("A1").Insert($xlShiftToRight)
$ws.Range("A1").Value = "part_of"

# --- GrossAnatomyAssociation: insert one column (part_of) ---
$ws = $wb.Worksheets.Item("GrossAnatomyAssociation")
$ws.Range("A1").Insert($xlShiftToRight)
$ws.Range("A1").Value = "part_of"
